$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 625485.5
$ws.Range("I9").Value = 708
$ws.Range("J9").Value = 1666781.4
$ws.Range("K9").Value = 708
$ws.Range("L9").Value = 1666781.4
$ws.Range("M9").Value = -539
$ws.Range("N9").Value = -1667119.4

$ws.Range("H33").Value = 448.7143
$ws.Range("I33").Value = 338.72726
$ws.Range("J33").Value = 852
$ws.Range("K33").Value = 338.72726
$ws.Range("L33").Value = 852
$ws.Range("M33").Value = -109.72726
$ws.Range("N33").Value = -1310

$ws.Range("H64").Value = 7402.6665
$ws.Range("I64").Value = 6105.5
$ws.Range("J64").Value = 9997
$ws.Range("K64").Value = 6105.5
$ws.Range("L64").Value = 9997
$ws.Range("M64").Value = -5857.5
$ws.Range("N64").Value = -10493

$ws.Range("H67").Value = 7402.6665
$ws.Range("I67").Value = 6105.5
$ws.Range("J67").Value = 9997
$ws.Range("K67").Value = 6105.5
$ws.Range("L67").Value = 9997
$ws.Range("M67").Value = -5247.5
$ws.Range("N67").Value = -11713

$ws.Range("H80").Value = 129.2
$ws.Range("I80").Value = 32.57143
$ws.Range("J80").Value = 213.75
$ws.Range("K80").Value = 97.71429000000001
$ws.Range("L80").Value = 641.25
$ws.Range("M80").Value = 900.28571
$ws.Range("N80").Value = -2637.25

$ws.Range("H83").Value = 129.2
$ws.Range("I83").Value = 32.57143
$ws.Range("J83").Value = 213.75
$ws.Range("K83").Value = 293.14287
$ws.Range("L83").Value = 1923.75
$ws.Range("M83").Value = 4698.85713
$ws.Range("N83").Value = -11907.75

$ws.Range("H100").Value = 998.4286
$ws.Range("I100").Value = 998.2
$ws.Range("J100").Value = 999
$ws.Range("K100").Value = 998.2
$ws.Range("L100").Value = 999
$ws.Range("M100").Value = -457.2
$ws.Range("N100").Value = -2081

$ws.Range("H106").Value = 1900
$ws.Range("I106").Value = 1900
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 1900
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1269

$ws.Range("H132").Value = 4180.857
$ws.Range("I132").Value = 3711.6667
$ws.Range("J132").Value = 6996
$ws.Range("K132").Value = 11135.0001
$ws.Range("L132").Value = 20988
$ws.Range("M132").Value = -8605.000100000001
$ws.Range("N132").Value = -26048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2248.5
$ws.Range("I2").Value = 2319.2727
$ws.Range("J2").Value = 1470
$ws.Range("K2").Value = 2319.2727
$ws.Range("L2").Value = 1470
$ws.Range("M2").Value = -2206.2727
$ws.Range("N2").Value = -1696

$ws.Range("H32").Value = 3314.8718
$ws.Range("I32").Value = 1924.8055
$ws.Range("J32").Value = 19995.666
$ws.Range("K32").Value = 1924.8055
$ws.Range("L32").Value = 19995.666
$ws.Range("M32").Value = -1637.8055
$ws.Range("N32").Value = -20569.666

$ws.Range("H74").Value = 1448.2084
$ws.Range("I74").Value = 1312.381
$ws.Range("J74").Value = 2399
$ws.Range("K74").Value = 1312.381
$ws.Range("L74").Value = 2399
$ws.Range("M74").Value = -438.3810000000001
$ws.Range("N74").Value = -4147

$ws.Range("H77").Value = 1448.2084
$ws.Range("I77").Value = 1312.381
$ws.Range("J77").Value = 2399
$ws.Range("K77").Value = 6561.905000000001
$ws.Range("L77").Value = 11995
$ws.Range("M77").Value = -2193.905000000001
$ws.Range("N77").Value = -20731

$ws.Range("H97").Value = 1654.9375
$ws.Range("I97").Value = 1540.7142
$ws.Range("J97").Value = 2454.5
$ws.Range("K97").Value = 1540.7142
$ws.Range("L97").Value = 2454.5
$ws.Range("M97").Value = -1044.7142
$ws.Range("N97").Value = -3446.5

$ws.Range("H102").Value = 2428.6667
$ws.Range("I102").Value = 2428.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2428.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -806.6667000000002

$ws.Range("H116").Value = 2248.5
$ws.Range("I116").Value = 2319.2727
$ws.Range("J116").Value = 1470
$ws.Range("K116").Value = 2319.2727
$ws.Range("L116").Value = 1470
$ws.Range("M116").Value = -25.27269999999999
$ws.Range("N116").Value = -6058

$ws.Range("H132").Value = 2932.762
$ws.Range("I132").Value = 2599.4375
$ws.Range("J132").Value = 3999.4
$ws.Range("K132").Value = 7798.3125
$ws.Range("L132").Value = 11998.2
$ws.Range("M132").Value = -5268.3125
$ws.Range("N132").Value = -17058.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2248.5
$ws.Range("I3").Value = 2319.2727
$ws.Range("J3").Value = 1470
$ws.Range("K3").Value = 2319.2727
$ws.Range("L3").Value = 1470
$ws.Range("M3").Value = -2205.2727
$ws.Range("N3").Value = -1698

$ws.Range("H80").Value = 1585.625
$ws.Range("I80").Value = 53.5
$ws.Range("J80").Value = 2096.3333
$ws.Range("K80").Value = 53.5
$ws.Range("L80").Value = 2096.3333
$ws.Range("M80").Value = 944.5
$ws.Range("N80").Value = -4092.3333

$ws.Range("H83").Value = 1585.625
$ws.Range("I83").Value = 53.5
$ws.Range("J83").Value = 2096.3333
$ws.Range("K83").Value = 267.5
$ws.Range("L83").Value = 10481.6665
$ws.Range("M83").Value = 4724.5
$ws.Range("N83").Value = -20465.6665

$ws.Range("H94").Value = 3238.6667
$ws.Range("I94").Value = 1007.3333
$ws.Range("J94").Value = 9932.666999999999
$ws.Range("K94").Value = 1007.3333
$ws.Range("L94").Value = 9932.666999999999
$ws.Range("M94").Value = -556.3333
$ws.Range("N94").Value = -10834.667

$ws.Range("H99").Value = 1955
$ws.Range("I99").Value = 1955
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1955
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -457

$ws.Range("H134").Value = 1783.2222
$ws.Range("I134").Value = 1783.2222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5349.6666
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2814.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3396.5
$ws.Range("I62").Value = 3396.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3396.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2772.5

$ws.Range("H65").Value = 3396.5
$ws.Range("I65").Value = 3396.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 16982.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -13862.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H68").Value = 3855.625
$ws.Range("I68").Value = 4290.6665
$ws.Range("J68").Value = 3594.6
$ws.Range("K68").Value = 12871.9995
$ws.Range("L68").Value = 10783.8
$ws.Range("M68").Value = -12060.9995
$ws.Range("N68").Value = -12405.8

$ws.Range("H71").Value = 3855.625
$ws.Range("I71").Value = 4290.6665
$ws.Range("J71").Value = 3594.6
$ws.Range("K71").Value = 38615.9985
$ws.Range("L71").Value = 32351.4
$ws.Range("M71").Value = -34559.9985
$ws.Range("N71").Value = -40463.39999999999

$ws.Range("H92").Value = 237.5
$ws.Range("I92").Value = 199
$ws.Range("J92").Value = 250.33333
$ws.Range("K92").Value = 597
$ws.Range("L92").Value = 750.99999
$ws.Range("M92").Value = 651
$ws.Range("N92").Value = -3246.99999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5224
$ws.Range("M5").ClearContents()

$ws.Range("H70").Value = 7408.647
$ws.Range("I70").Value = 5326.1665
$ws.Range("J70").Value = 8544.546
$ws.Range("K70").Value = 5326.1665
$ws.Range("L70").Value = 8544.546
$ws.Range("M70").Value = -5056.1665
$ws.Range("N70").Value = -9084.546

$ws.Range("H73").Value = 7408.647
$ws.Range("I73").Value = 5326.1665
$ws.Range("J73").Value = 8544.546
$ws.Range("K73").Value = 5326.1665
$ws.Range("L73").Value = 8544.546
$ws.Range("M73").Value = -4390.1665
$ws.Range("N73").Value = -10416.546

$ws.Range("H122").Value = 1514.1111
$ws.Range("I122").Value = 1228.375
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 3685.125
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -1235.125
$ws.Range("N122").Value = -16300

$ws.Range("H126").Value = 3999.2
$ws.Range("I126").Value = 3425.2
$ws.Range("J126").Value = 4573.2
$ws.Range("K126").Value = 10275.6
$ws.Range("L126").Value = 13719.6
$ws.Range("M126").Value = -7805.599999999999
$ws.Range("N126").Value = -18659.6

$ws.Range("H132").Value = 3494.2273
$ws.Range("I132").Value = 2940.842
$ws.Range("J132").Value = 6999
$ws.Range("K132").Value = 8822.526
$ws.Range("L132").Value = 20997
$ws.Range("M132").Value = -6292.526
$ws.Range("N132").Value = -26057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 10000
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224

$ws.Range("H22").Value = 1868.2
$ws.Range("I22").Value = 498.5
$ws.Range("J22").Value = 2781.3333
$ws.Range("K22").Value = 498.5
$ws.Range("L22").Value = 2781.3333
$ws.Range("M22").Value = -203.5
$ws.Range("N22").Value = -3371.3333

$ws.Range("H27").Value = 1868.2
$ws.Range("I27").Value = 498.5
$ws.Range("J27").Value = 2781.3333
$ws.Range("K27").Value = 498.5
$ws.Range("L27").Value = 2781.3333
$ws.Range("M27").Value = -391.5
$ws.Range("N27").Value = -2995.3333

$ws.Range("H40").Value = 6640.8335
$ws.Range("I40").Value = 6282.8335
$ws.Range("J40").Value = 6998.8335
$ws.Range("K40").Value = 6282.8335
$ws.Range("L40").Value = 6998.8335
$ws.Range("M40").Value = -6146.8335
$ws.Range("N40").Value = -7270.8335

$ws.Range("H46").Value = 1598.9286
$ws.Range("I46").Value = 771.4286
$ws.Range("J46").Value = 2426.4285
$ws.Range("K46").Value = 771.4286
$ws.Range("L46").Value = 2426.4285
$ws.Range("M46").Value = -583.4286
$ws.Range("N46").Value = -2802.4285

$ws.Range("H68").Value = 2643.4285
$ws.Range("I68").Value = 2300.8
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 2300.8
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -1551.8
$ws.Range("N68").Value = -4998

$ws.Range("H71").Value = 2643.4285
$ws.Range("I71").Value = 2300.8
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 11504
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -7760
$ws.Range("N71").Value = -24988

$ws.Range("H100").Value = 2730.5386
$ws.Range("I100").Value = 2416.1667
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2416.1667
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1875.1667
$ws.Range("N100").Value = -4082

$ws.Range("H108").Value = 73749.75
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 73749.75
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 73749.75
$ws.Range("N108").Value = -81429.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 19500
$ws.Range("I2").Value = 19500
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 19500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -19388
$ws.Range("N2").ClearContents()

$ws.Range("H41").Value = 13883.5
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 13883.5
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 13883.5
$ws.Range("N41").Value = -14663.5

$ws.Range("H107").Value = 763.16
$ws.Range("I107").Value = 776.75
$ws.Range("J107").Value = 750.61536
$ws.Range("K107").Value = 2330.25
$ws.Range("L107").Value = 2251.84608
$ws.Range("M107").Value = -410.25
$ws.Range("N107").Value = -6091.84608

$ws.Range("H136").Value = 1455.375
$ws.Range("I136").Value = 1152.4
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 3577.5999
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -907.2000000000003
$ws.Range("N136").Value = -23100
